# Generate Report for Handback
# - Updates the zh-cn / de-de handback status text on the Overview sheet
# - Populates the "Latest Target File", "Latest Handback File" and
#   "Latest Handback DateTime" columns on the per-locale sheets now that the
#   handback has completed
# - Adds a hyperlink on the newly-populated "Latest Target File" cell

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$mdFileName = "0a24f56a-195f-4297-a987-da646670bcb0.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1fecb61ca8a0badef573aba574e902cacc70e0e4/e2e/0a24f56a-195f-4297-a987-da646670bcb0.md"

# ---- Overview sheet: both locale status cells move from "Ready for
# ---- handoff" to the new in-sync status text.
$overview = $wb.Worksheets.Item("Overview")
$overview.Cells.Item(2, 5).Value = $statusText   # E2 (zh-cn)
$overview.Cells.Item(2, 6).Value = $statusText   # F2 (de-de)
$overview.Columns.Item(5).EntireColumn.AutoFit()
$overview.Columns.Item(6).EntireColumn.AutoFit()

# ---- zh-cn sheet
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Cells.Item(2, 3).Value = $statusText        # C2 Status
$zhcn.Cells.Item(2, 9).Value = $mdFileName        # I2 Latest Target File
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrl, "", "", $mdFileName) | Out-Null
$zhcn.Cells.Item(2, 10).Value = "0a24f56a-195f-4297-a987-da646670bcb0.d72723ae0fc48b1b187dda5e4a8202c5d149647f.zh-cn.xlf"  # J2 Latest Handback File
$zhcn.Cells.Item(2, 11).Value = "2016-09-01 19:09:18"  # K2 Latest Handback DateTime
$zhcn.Columns.Item(3).EntireColumn.AutoFit()
$zhcn.Columns.Item(9).EntireColumn.AutoFit()
$zhcn.Columns.Item(10).EntireColumn.AutoFit()

# ---- de-de sheet
$dede = $wb.Worksheets.Item("de-de")
$dede.Cells.Item(2, 3).Value = $statusText        # C2 Status
$dede.Cells.Item(2, 9).Value = $mdFileName        # I2 Latest Target File
$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrl, "", "", $mdFileName) | Out-Null
$dede.Cells.Item(2, 10).Value = "0a24f56a-195f-4297-a987-da646670bcb0.d72723ae0fc48b1b187dda5e4a8202c5d149647f.de-de.xlf"  # J2 Latest Handback File
$dede.Cells.Item(2, 11).Value = "2016-09-01 19:09:26"  # K2 Latest Handback DateTime
$dede.Columns.Item(3).EntireColumn.AutoFit()
$dede.Columns.Item(9).EntireColumn.AutoFit()
$dede.Columns.Item(10).EntireColumn.AutoFit()
